$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 521.5
$ws.Cells.Item(12, 10).Value = 1029.2
$ws.Cells.Item(12, 12).Value = 1029.2
$ws.Cells.Item(12, 14).Value = -1369.2
$ws.Cells.Item(64, 8).Value = 7476.7
$ws.Cells.Item(64, 9).Value = 5695
$ws.Cells.Item(64, 10).Value = 7922.125
$ws.Cells.Item(64, 11).Value = 5695
$ws.Cells.Item(64, 12).Value = 7922.125
$ws.Cells.Item(64, 13).Value = -5447
$ws.Cells.Item(64, 14).Value = -8418.125
$ws.Cells.Item(67, 8).Value = 7476.7
$ws.Cells.Item(67, 9).Value = 5695
$ws.Cells.Item(67, 10).Value = 7922.125
$ws.Cells.Item(67, 11).Value = 5695
$ws.Cells.Item(67, 12).Value = 7922.125
$ws.Cells.Item(67, 13).Value = -4837
$ws.Cells.Item(67, 14).Value = -9638.125
$ws.Cells.Item(86, 8).Value = 2510.8462
$ws.Cells.Item(86, 9).Value = 2946.1667
$ws.Cells.Item(86, 11).Value = 2946.1667
$ws.Cells.Item(86, 13).Value = -1823.1667
$ws.Cells.Item(89, 8).Value = 2510.8462
$ws.Cells.Item(89, 9).Value = 2946.1667
$ws.Cells.Item(89, 11).Value = 14730.8335
$ws.Cells.Item(89, 13).Value = -9114.833500000001
$ws.Cells.Item(134, 8).Value = 32994.688
$ws.Cells.Item(134, 10).Value = 32994.688
$ws.Cells.Item(134, 12).Value = 32994.688
$ws.Cells.Item(134, 14).Value = -43134.688

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1447.625
$ws.Cells.Item(2, 9).Value = 793.3333
$ws.Cells.Item(2, 10).Value = 1840.2
$ws.Cells.Item(2, 11).Value = 793.3333
$ws.Cells.Item(2, 12).Value = 1840.2
$ws.Cells.Item(2, 13).Value = -680.3333
$ws.Cells.Item(2, 14).Value = -2066.2
$ws.Cells.Item(32, 8).Value = 2554.5
$ws.Cells.Item(32, 9).Value = 2336.158
$ws.Cells.Item(32, 11).Value = 2336.158
$ws.Cells.Item(32, 13).Value = -2049.158
$ws.Cells.Item(116, 8).Value = 1447.625
$ws.Cells.Item(116, 9).Value = 793.3333
$ws.Cells.Item(116, 10).Value = 1840.2
$ws.Cells.Item(116, 11).Value = 793.3333
$ws.Cells.Item(116, 12).Value = 1840.2
$ws.Cells.Item(116, 13).Value = 1500.6667
$ws.Cells.Item(116, 14).Value = -6428.2
$ws.Cells.Item(132, 8).Value = 3230.8408
$ws.Cells.Item(132, 9).Value = 2903.5881
$ws.Cells.Item(132, 10).Value = 4343.5
$ws.Cells.Item(132, 11).Value = 8710.764299999999
$ws.Cells.Item(132, 12).Value = 13030.5
$ws.Cells.Item(132, 13).Value = -6180.764299999999
$ws.Cells.Item(132, 14).Value = -18090.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1447.625
$ws.Cells.Item(3, 9).Value = 793.3333
$ws.Cells.Item(3, 10).Value = 1840.2
$ws.Cells.Item(3, 11).Value = 793.3333
$ws.Cells.Item(3, 12).Value = 1840.2
$ws.Cells.Item(3, 13).Value = -679.3333
$ws.Cells.Item(3, 14).Value = -2068.2
$ws.Cells.Item(8, 8).Value = 550
$ws.Cells.Item(8, 9).Value = 500
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 11).Value = 500
$ws.Cells.Item(8, 12).Value = 600
$ws.Cells.Item(8, 13).Value = -360
$ws.Cells.Item(8, 14).Value = -880
$ws.Cells.Item(14, 8).Value = 1500
$ws.Cells.Item(14, 10).Value = 1500
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 14).Value = -1844
$ws.Cells.Item(17, 8).Value = 12999
$ws.Cells.Item(17, 10).Value = 12999
$ws.Cells.Item(17, 12).Value = 12999
$ws.Cells.Item(17, 14).Value = -13343
$ws.Cells.Item(132, 8).Value = 89999.5
$ws.Cells.Item(132, 10).Value = 89999.5
$ws.Cells.Item(132, 12).Value = 89999.5
$ws.Cells.Item(132, 14).Value = -100119.5
$ws.Cells.Item(133, 8).Value = 53994.75
$ws.Cells.Item(133, 9).Value = 35000
$ws.Cells.Item(133, 10).Value = 60326.332
$ws.Cells.Item(133, 11).Value = 35000
$ws.Cells.Item(133, 12).Value = 60326.332
$ws.Cells.Item(133, 13).Value = -29940
$ws.Cells.Item(133, 14).Value = -70446.33199999999
$ws.Cells.Item(139, 8).Value = 162617.67
$ws.Cells.Item(139, 10).Value = 184999.4
$ws.Cells.Item(139, 12).Value = 184999.4
$ws.Cells.Item(139, 14).Value = -195279.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 1851.5
$ws.Cells.Item(12, 10).Value = 2467
$ws.Cells.Item(12, 12).Value = 2467
$ws.Cells.Item(12, 14).Value = -2807
$ws.Cells.Item(14, 8).Value = 2155
$ws.Cells.Item(14, 9).Value = 2810
$ws.Cells.Item(14, 10).Value = 1500
$ws.Cells.Item(14, 11).Value = 2810
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = -2640
$ws.Cells.Item(14, 14).Value = -1840
$ws.Cells.Item(31, 8).Value = 3663.0789
$ws.Cells.Item(31, 9).Value = 2236
$ws.Cells.Item(31, 11).Value = 2236
$ws.Cells.Item(31, 13).Value = -1941
$ws.Cells.Item(34, 8).Value = 3663.0789
$ws.Cells.Item(34, 9).Value = 2236
$ws.Cells.Item(34, 11).Value = 2236
$ws.Cells.Item(34, 13).Value = -2034
$ws.Cells.Item(62, 8).Value = 7275.1113
$ws.Cells.Item(62, 9).Value = 6955.6
$ws.Cells.Item(62, 10).Value = 7674.5
$ws.Cells.Item(62, 11).Value = 6955.6
$ws.Cells.Item(62, 12).Value = 7674.5
$ws.Cells.Item(62, 13).Value = -6331.6
$ws.Cells.Item(62, 14).Value = -8922.5
$ws.Cells.Item(65, 8).Value = 7275.1113
$ws.Cells.Item(65, 9).Value = 6955.6
$ws.Cells.Item(65, 10).Value = 7674.5
$ws.Cells.Item(65, 11).Value = 34778
$ws.Cells.Item(65, 12).Value = 38372.5
$ws.Cells.Item(65, 13).Value = -31658
$ws.Cells.Item(65, 14).Value = -44612.5
$ws.Cells.Item(105, 8).Value = 1088.5834
$ws.Cells.Item(105, 9).Value = 1004.7273
$ws.Cells.Item(105, 11).Value = 1004.7273
$ws.Cells.Item(105, 13).Value = 742.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 112.73913
$ws.Cells.Item(2, 10).Value = 102.3
$ws.Cells.Item(2, 12).Value = 613.8
$ws.Cells.Item(2, 14).Value = -839.8
$ws.Cells.Item(5, 8).Value = 896.875
$ws.Cells.Item(5, 9).Value = 292
$ws.Cells.Item(5, 11).Value = 876
$ws.Cells.Item(5, 13).Value = -764
$ws.Cells.Item(23, 8).Value = 516.5
$ws.Cells.Item(23, 9).Value = 295.16666
$ws.Cells.Item(23, 10).Value = 611.3570999999999
$ws.Cells.Item(23, 11).Value = 885.4999799999999
$ws.Cells.Item(23, 12).Value = 1834.0713
$ws.Cells.Item(23, 13).Value = -650.4999799999999
$ws.Cells.Item(23, 14).Value = -2304.0713
$ws.Cells.Item(135, 8).Value = 896.875
$ws.Cells.Item(135, 9).Value = 292
$ws.Cells.Item(135, 11).Value = 2628
$ws.Cells.Item(135, 13).Value = -93
$ws.Cells.Item(141, 8).Value = 22327.6
$ws.Cells.Item(141, 9).Value = 4947.5835
$ws.Cells.Item(141, 11).Value = 14842.7505
$ws.Cells.Item(141, 13).Value = -9662.750499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1868.6923
$ws.Cells.Item(80, 9).Value = 1813.5714
$ws.Cells.Item(80, 10).Value = 1933
$ws.Cells.Item(80, 11).Value = 1813.5714
$ws.Cells.Item(80, 12).Value = 1933
$ws.Cells.Item(80, 13).Value = -815.5714
$ws.Cells.Item(80, 14).Value = -3929
$ws.Cells.Item(83, 8).Value = 1868.6923
$ws.Cells.Item(83, 9).Value = 1813.5714
$ws.Cells.Item(83, 10).Value = 1933
$ws.Cells.Item(83, 11).Value = 9067.857
$ws.Cells.Item(83, 12).Value = 9665
$ws.Cells.Item(83, 13).Value = -4075.857
$ws.Cells.Item(83, 14).Value = -19649
$ws.Cells.Item(107, 8).Value = 461.21738
$ws.Cells.Item(107, 9).Value = 250.28572
$ws.Cells.Item(107, 10).Value = 789.3333
$ws.Cells.Item(107, 11).Value = 250.28572
$ws.Cells.Item(107, 12).Value = 789.3333
$ws.Cells.Item(107, 13).Value = 1669.71428
$ws.Cells.Item(107, 14).Value = -4629.3333
$ws.Cells.Item(122, 8).Value = 2347.6667
$ws.Cells.Item(122, 9).Value = 2242.9546
$ws.Cells.Item(122, 10).Value = 2808.4
$ws.Cells.Item(122, 11).Value = 6728.8638
$ws.Cells.Item(122, 12).Value = 8425.200000000001
$ws.Cells.Item(122, 13).Value = -4278.8638
$ws.Cells.Item(122, 14).Value = -13325.2
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2196.96
$ws.Cells.Item(132, 9).Value = 1900.3125
$ws.Cells.Item(132, 11).Value = 5700.9375
$ws.Cells.Item(132, 13).Value = -3170.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(31, 8).Value = 9014.111000000001
$ws.Cells.Item(31, 9).Value = 1500
$ws.Cells.Item(31, 10).Value = 9953.375
$ws.Cells.Item(31, 11).Value = 1500
$ws.Cells.Item(31, 12).Value = 9953.375
$ws.Cells.Item(31, 13).Value = -1252
$ws.Cells.Item(31, 14).Value = -10449.375
$ws.Cells.Item(68, 8).Value = 2101.32
$ws.Cells.Item(68, 9).Value = 1932.3889
$ws.Cells.Item(68, 11).Value = 1932.3889
$ws.Cells.Item(68, 13).Value = -1183.3889
$ws.Cells.Item(71, 8).Value = 2101.32
$ws.Cells.Item(71, 9).Value = 1932.3889
$ws.Cells.Item(71, 11).Value = 9661.9445
$ws.Cells.Item(71, 13).Value = -5917.9445
$ws.Cells.Item(82, 8).Value = 2736.9033
$ws.Cells.Item(82, 10).Value = 7705.857
$ws.Cells.Item(82, 12).Value = 7705.857
$ws.Cells.Item(82, 14).Value = -8427.857
$ws.Cells.Item(85, 8).Value = 2736.9033
$ws.Cells.Item(85, 10).Value = 7705.857
$ws.Cells.Item(85, 12).Value = 7705.857
$ws.Cells.Item(85, 14).Value = -10201.857
$ws.Cells.Item(93, 8).Value = 1114.5385
$ws.Cells.Item(93, 9).Value = 1064.5555
$ws.Cells.Item(93, 10).Value = 1227
$ws.Cells.Item(93, 11).Value = 1064.5555
$ws.Cells.Item(93, 12).Value = 1227
$ws.Cells.Item(93, 13).Value = 183.4445000000001
$ws.Cells.Item(93, 14).Value = -3723
$ws.Cells.Item(136, 8).Value = 1548.5
$ws.Cells.Item(136, 9).Value = 1120.2667
$ws.Cells.Item(136, 11).Value = 3360.800099999999
$ws.Cells.Item(136, 13).Value = -810.8000999999995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 6357
$ws.Cells.Item(81, 9).Value = 6999.909
$ws.Cells.Item(81, 11).Value = 13999.818
$ws.Cells.Item(81, 13).Value = -12938.818
$ws.Cells.Item(84, 8).Value = 6357
$ws.Cells.Item(84, 9).Value = 6999.909
$ws.Cells.Item(84, 11).Value = 69999.09
$ws.Cells.Item(84, 13).Value = -64695.09
$ws.Cells.Item(126, 8).Value = 2322
$ws.Cells.Item(126, 9).Value = 2000.5
$ws.Cells.Item(126, 10).Value = 2579.2
$ws.Cells.Item(126, 11).Value = 6001.5
$ws.Cells.Item(126, 12).Value = 7737.599999999999
$ws.Cells.Item(126, 13).Value = -3531.5
$ws.Cells.Item(126, 14).Value = -12677.6
